$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.990.84'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '1.858.08'
$ws.Range("E3").Value = '  -0.91%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '312.45'
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").Value = '0.5133'
$ws.Range("E7").Value = '  +0.99%  '
$ws.Range("D8").Value = '0.3830'
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("D9").Value = '0.08239'
$ws.Range("E9").Value = '  -8.68%  '
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("D11").Value = '41.49'
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").Value = '6.184'
$ws.Range("E12").Value = '  -2.54%  '
$ws.Range("D13").Value = '20.54'
$ws.Range("E13").Value = '  -0.91%  '
$ws.Range("D14").Value = '1.860.07'
$ws.Range("E14").Value = '  -0.21%  '
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("D16").Value = '1.003'
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("D18").Value = '90.57'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("D20").Value = '17.67'
$ws.Range("E20").Value = '  -2.75%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  -1.91%  '
$ws.Range("D23").Value = '28.016.50'
$ws.Range("E23").Value = '  -0.27%  '
$ws.Range("D24").Value = '11.06'
$ws.Range("E24").Value = '  -3.14%  '
$ws.Range("D25").Value = '2.244'
$ws.Range("E25").Value = '  -1.63%  '
$ws.Range("D26").Value = '2.073.12'
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("D27").Value = '2.508'
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("D28").Value = '158.05'
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("D29").Value = '20.44'
$ws.Range("E29").Value = '  -1.68%  '
$ws.Range("D30").Value = '124.56'
$ws.Range("E30").Value = '  -1.62%  '
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("E32").Value = '  -2.93%  '
$ws.Range("D33").Value = '5.943'
$ws.Range("E33").Value = '  +5.79%  '
$ws.Range("D34").Value = '3.595'
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D35").Value = '9.339'
$ws.Range("E35").Value = '  -3.27%  '
$ws.Range("D36").Value = '0.02414'
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("D37").Value = '0.06483'
$ws.Range("E37").Value = '  -1.45%  '
$ws.Range("D38").Value = '0.2167'
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("D39").Value = '0.6526'
$ws.Range("E39").Value = '  +2.06%  '
$ws.Range("D40").Value = '1.193'
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("D41").Value = '5.004'
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("D42").Value = '1.220'
$ws.Range("E42").Value = '  -3.59%  '
$ws.Range("D43").Value = '11.15'
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("D44").Value = '0.6143'
$ws.Range("E44").Value = '  +2.09%  '
$ws.Range("D45").Value = '13.01'
$ws.Range("E45").Value = '  -1.23%  '
$ws.Range("E46").Value = '  +0.32%  '
$ws.Range("D47").Value = '3.672'
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("E48").Value = '  +0.56%  '
$ws.Range("D50").Value = '120.34'
$ws.Range("E50").Value = '  -0.92%  '
$ws.Range("D51").Value = '78.29'
$ws.Range("E51").Value = '  -1.85%  '
